$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 04.02.2022 16:00"

# D2: change from text "+0.4" to numeric 0.4
$ws.Range("D2").Value = 0.4

# E2: change from text date string to numeric date serial,
# using the same number format style as the other date cells (e.g. E3)
$ws.Range("E2").Value = 44596.65634259259
$ws.Range("E2").NumberFormat = $ws.Range("E3").NumberFormat
